$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "reyan"
$ws.Range("A1").Value = "Vedhika Menon"
